$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        if ($cell.Value2 -eq "Ready for handoff") {
            $cell.Value = "In Translation"
        }
    }
}

$overview = $wb.Worksheets.Item("Overview")
$overview.Columns.Item(5).ColumnWidth = 13.4101845877511
$overview.Columns.Item(6).ColumnWidth = 13.4101845877511

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Columns.Item(3).ColumnWidth = 13.4101845877511

$dede = $wb.Worksheets.Item("de-de")
$dede.Columns.Item(3).ColumnWidth = 13.4101845877511
